$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are written as text, matching the source data
# (Coinranking prices are strings, sometimes using "." as a thousands separator,
# so they must not be auto-converted to numbers by Excel).

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '68.213.13'
$ws.Cells.Item(2, 5).Value = '  +3.71%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.634.65'
$ws.Cells.Item(3, 5).Value = '  +3.14%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.71%  '

$ws.Cells.Item(5, 5).Value = '  +10.90%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '569.57'
$ws.Cells.Item(6, 5).Value = '  -1.26%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.612.32'
$ws.Cells.Item(7, 5).Value = '  +2.71%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.617'
$ws.Cells.Item(8, 5).Value = '  +2.67%  '

$ws.Cells.Item(9, 5).Value = '  +0.01%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.677'
$ws.Cells.Item(10, 5).Value = '  +2.53%  '

$ws.Cells.Item(11, 2).Value = 'Avalanche'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '57.63'
$ws.Cells.Item(11, 5).Value = '  +7.66%  '

$ws.Cells.Item(12, 2).Value = 'Dogecoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.154'
$ws.Cells.Item(12, 5).Value = '  +8.94%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000291'
$ws.Cells.Item(13, 5).Value = '  +17.80%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '10.05'
$ws.Cells.Item(14, 5).Value = '  +3.81%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.223.00'
$ws.Cells.Item(15, 5).Value = '  +3.58%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.637.36'
$ws.Cells.Item(16, 5).Value = '  +3.91%  '

$ws.Cells.Item(17, 5).Value = '  +0.80%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.47'
$ws.Cells.Item(18, 5).Value = '  +4.00%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '68.168.90'
$ws.Cells.Item(19, 5).Value = '  +4.30%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.56'
$ws.Cells.Item(20, 5).Value = '  +2.65%  '

$ws.Cells.Item(21, 5).Value = '  +4.21%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '401.27'
$ws.Cells.Item(22, 5).Value = '  +3.27%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '12.92'
$ws.Cells.Item(23, 5).Value = '  +26.71%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '4.19'
$ws.Cells.Item(24, 5).Value = '  -1.38%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '85.70'
$ws.Cells.Item(25, 5).Value = '  +2.42%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.95'
$ws.Cells.Item(26, 5).Value = '  +4.27%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.59'
$ws.Cells.Item(27, 5).Value = '  +3.33%  '

$ws.Cells.Item(28, 5).Value = '  +1.97%  '

$ws.Cells.Item(29, 5).Value = '  +8.62%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '8.19'
$ws.Cells.Item(30, 5).Value = '  +22.47%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '9.14'
$ws.Cells.Item(31, 5).Value = '  +3.90%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '31.83'
$ws.Cells.Item(32, 5).Value = '  +4.10%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '702.07'
$ws.Cells.Item(33, 5).Value = '  +15.65%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '12.23'
$ws.Cells.Item(34, 5).Value = '  +2.98%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.116'
$ws.Cells.Item(35, 5).Value = '  +4.91%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '64.18'
$ws.Cells.Item(36, 5).Value = '  -1.19%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '42.65'
$ws.Cells.Item(37, 5).Value = '  +3.93%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.425'
$ws.Cells.Item(38, 5).Value = '  +15.52%  '

$ws.Cells.Item(39, 5).Value = '  +0.04%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0₃0781'
$ws.Cells.Item(40, 5).Value = '  +6.48%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.140'
$ws.Cells.Item(41, 5).Value = '  +8.62%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.229.32'
$ws.Cells.Item(42, 5).Value = '  +13.17%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.13'
$ws.Cells.Item(43, 5).Value = '  +14.29%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.80'
$ws.Cells.Item(44, 5).Value = '  +17.07%  '

$ws.Cells.Item(45, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  +0.68%  '

$ws.Cells.Item(46, 2).Value = 'dogwifhat'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.00'
$ws.Cells.Item(46, 5).Value = '  +38.87%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0420'
$ws.Cells.Item(47, 5).Value = '  +4.61%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.71'
$ws.Cells.Item(48, 5).Value = '  +11.55%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.91'
$ws.Cells.Item(49, 5).Value = '  +9.51%  '

$ws.Cells.Item(50, 5).Value = '  +2.26%  '

$ws.Cells.Item(51, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '3.09'
$ws.Cells.Item(51, 5).Value = '  +4.90%  '
